$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data edits -------------------------------------------------------
# Row 4: "perkinazo" -> "jorsi"
$ws.Range("B4").Value = "jorsi"

# Row 6 ("nelson" / "solicitud working" row) is removed - contents cleared.
$ws.Range("A6:F6").ClearContents() | Out-Null

# --- Emails as hyperlinks (replacing the placeholder strings) -----------
$ws.Range("E2").Value = "thomas@redsalud.dirsal"
$ws.Hyperlinks.Add($ws.Range("E2"), "mailto:thomas@redsalud.dirsal") | Out-Null

$ws.Range("E3").Value = "leandro@redsalud.dirsal"
$ws.Hyperlinks.Add($ws.Range("E3"), "mailto:leandro@redsalud.dirsal") | Out-Null

$ws.Range("E4").Value = "jorsi@redsalud.dirsal"
$ws.Hyperlinks.Add($ws.Range("E4"), "mailto:jorsi@redsalud.dirsal") | Out-Null

# Row 5: old "jorsi" duplicate row is repurposed to hold the data that
# used to live in row 6 ("nelson" / rut 12345678), renamed to "pretoriano"
# and switched to "activo".
$ws.Range("A5").Value = 12345678
$ws.Range("B5").Value = "pretoriano"
$ws.Range("C5").Value = "activo"

$ws.Range("E5").Value = "pretorian@mail"
$ws.Hyperlinks.Add($ws.Range("E5"), "mailto:pretorian@mail") | Out-Null

# --- Formatting ---------------------------------------------------------
# The new row 5 / blank row 6 get an underline format applied (applied
# after the hyperlinks, so the hyperlink style wins on E5).
$ws.Range("A5:F6").Font.Underline = $true

# --- Selection ------------------------------------------------------------
$ws.Range("F6").Select() | Out-Null
